$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 428.2857
$ws.Range("J4").Value = 724.5
$ws.Range("L4").Value = 724.5
$ws.Range("N4").Value = -952.5
$ws.Range("H19").Value = 980.7083
$ws.Range("I19").Value = 943.2857
$ws.Range("J19").Value = 1033.1
$ws.Range("K19").Value = 943.2857
$ws.Range("L19").Value = 1033.1
$ws.Range("M19").Value = -768.2857
$ws.Range("N19").Value = -1383.1
$ws.Range("H58").Value = 6209.9
$ws.Range("I58").Value = 1366.3334
$ws.Range("J58").Value = 8285.714
$ws.Range("K58").Value = 4099.0002
$ws.Range("L58").Value = 24857.142
$ws.Range("M58").Value = -3949.0002
$ws.Range("N58").Value = -25157.142
$ws.Range("H76").Value = 6699
$ws.Range("I76").Value = 5666.8423
$ws.Range("J76").Value = 9500.571
$ws.Range("K76").Value = 5666.8423
$ws.Range("L76").Value = 9500.571
$ws.Range("M76").Value = -5351.8423
$ws.Range("N76").Value = -10130.571
$ws.Range("H79").Value = 6699
$ws.Range("I79").Value = 5666.8423
$ws.Range("J79").Value = 9500.571
$ws.Range("K79").Value = 5666.8423
$ws.Range("L79").Value = 9500.571
$ws.Range("M79").Value = -4574.8423
$ws.Range("N79").Value = -11684.571
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("H100").Value = 2194.7144
$ws.Range("I100").Value = 613.2
$ws.Range("K100").Value = 613.2
$ws.Range("M100").Value = -72.20000000000005
$ws.Range("H132").Value = 1633.8594
$ws.Range("I132").Value = 1525.2742
$ws.Range("K132").Value = 4575.8226
$ws.Range("M132").Value = -2045.8226
$ws.Range("H138").Value = 5675.8535
$ws.Range("I138").Value = 2848.2
$ws.Range("J138").Value = 6110.877
$ws.Range("K138").Value = 8544.599999999999
$ws.Range("L138").Value = 18332.631
$ws.Range("M138").Value = -3404.599999999999
$ws.Range("N138").Value = -28612.631

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6463.4736
$ws.Range("I2").Value = 1178.2307
$ws.Range("K2").Value = 1178.2307
$ws.Range("M2").Value = -1065.2307
$ws.Range("H32").Value = 10931.125
$ws.Range("I32").Value = 8775.772999999999
$ws.Range("K32").Value = 8775.772999999999
$ws.Range("M32").Value = -8488.772999999999
$ws.Range("H41").Value = 2781.1428
$ws.Range("I41").Value = 1493.6
$ws.Range("K41").Value = 1493.6
$ws.Range("M41").Value = -1079.6
$ws.Range("H61").Value = 4583.1353
$ws.Range("I61").Value = 3860.2307
$ws.Range("J61").Value = 6291.8184
$ws.Range("K61").Value = 3860.2307
$ws.Range("L61").Value = 6291.8184
$ws.Range("M61").Value = -3648.2307
$ws.Range("N61").Value = -6715.8184
$ws.Range("H74").Value = 2948.96
$ws.Range("I74").Value = 2290.4666
$ws.Range("K74").Value = 2290.4666
$ws.Range("M74").Value = -1416.4666
$ws.Range("H77").Value = 2948.96
$ws.Range("I77").Value = 2290.4666
$ws.Range("K77").Value = 11452.333
$ws.Range("M77").Value = -7084.333000000001
$ws.Range("H102").Value = 2757.9375
$ws.Range("I102").Value = 2779
$ws.Range("K102").Value = 2779
$ws.Range("M102").Value = -1157
$ws.Range("H110").Value = 3888.8215
$ws.Range("I110").Value = 3357.6086
$ws.Range("J110").Value = 6332.4
$ws.Range("K110").Value = 3357.6086
$ws.Range("L110").Value = 6332.4
$ws.Range("M110").Value = -1312.6086
$ws.Range("N110").Value = -10422.4
$ws.Range("H116").Value = 6463.4736
$ws.Range("I116").Value = 1178.2307
$ws.Range("K116").Value = 1178.2307
$ws.Range("M116").Value = 1115.7693
$ws.Range("H122").Value = 4414.6816
$ws.Range("I122").Value = 3985
$ws.Range("K122").Value = 11955
$ws.Range("M122").Value = -9505
$ws.Range("H132").Value = 3323.7856
$ws.Range("I132").Value = 2962.4722
$ws.Range("J132").Value = 5491.6665
$ws.Range("K132").Value = 8887.4166
$ws.Range("L132").Value = 16474.9995
$ws.Range("M132").Value = -6357.4166
$ws.Range("N132").Value = -21534.9995
$ws.Range("H136").Value = 4583.1353
$ws.Range("I136").Value = 3860.2307
$ws.Range("J136").Value = 6291.8184
$ws.Range("K136").Value = 11580.6921
$ws.Range("L136").Value = 18875.4552
$ws.Range("M136").Value = -9030.6921
$ws.Range("N136").Value = -23975.4552

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6463.4736
$ws.Range("I3").Value = 1178.2307
$ws.Range("K3").Value = 1178.2307
$ws.Range("M3").Value = -1064.2307
$ws.Range("H64").Value = 1410.7142
$ws.Range("J64").Value = 1410.7142
$ws.Range("L64").Value = 1410.7142
$ws.Range("N64").Value = -1860.7142
$ws.Range("H67").Value = 1410.7142
$ws.Range("J67").Value = 1410.7142
$ws.Range("L67").Value = 1410.7142
$ws.Range("N67").Value = -2970.7142
$ws.Range("H94").Value = 2406.3845
$ws.Range("I94").Value = 2382.64
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 2382.64
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -1931.64
$ws.Range("N94").Value = -3902
$ws.Range("H99").Value = 3731.875
$ws.Range("I99").Value = 1852.5454
$ws.Range("K99").Value = 1852.5454
$ws.Range("M99").Value = -354.5454
$ws.Range("H134").Value = 4665.7896
$ws.Range("I134").Value = 4322.467
$ws.Range("J134").Value = 5953.25
$ws.Range("K134").Value = 12967.401
$ws.Range("L134").Value = 17859.75
$ws.Range("M134").Value = -10432.401
$ws.Range("N134").Value = -22929.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38232
$ws.Range("I31").Value = 3358.25
$ws.Range("J31").Value = 50362
$ws.Range("K31").Value = 3358.25
$ws.Range("L31").Value = 50362
$ws.Range("M31").Value = -3063.25
$ws.Range("N31").Value = -50952
$ws.Range("H34").Value = 38232
$ws.Range("I34").Value = 3358.25
$ws.Range("J34").Value = 50362
$ws.Range("K34").Value = 3358.25
$ws.Range("L34").Value = 50362
$ws.Range("M34").Value = -3156.25
$ws.Range("N34").Value = -50766
$ws.Range("H105").Value = 3074
$ws.Range("I105").Value = 1182.8334
$ws.Range("J105").Value = 4965.1665
$ws.Range("K105").Value = 1182.8334
$ws.Range("L105").Value = 4965.1665
$ws.Range("M105").Value = 564.1666
$ws.Range("N105").Value = -8459.166499999999
$ws.Range("H122").Value = 3488.5134
$ws.Range("I122").Value = 1806.4814
$ws.Range("J122").Value = 8030
$ws.Range("K122").Value = 5419.4442
$ws.Range("L122").Value = 24090
$ws.Range("M122").Value = -2969.4442
$ws.Range("N122").Value = -28990
$ws.Range("H134").Value = 2754.0454
$ws.Range("I134").Value = 2135.5264
$ws.Range("K134").Value = 6406.5792
$ws.Range("M134").Value = -3871.5792

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3860
$ws.Range("J113").Value = 3860
$ws.Range("L113").Value = 11580
$ws.Range("N113").Value = -15920
$ws.Range("H122").Value = 3282.125
$ws.Range("J122").Value = 3703.55
$ws.Range("L122").Value = 33331.95
$ws.Range("N122").Value = -38231.95
$ws.Range("H137").Value = 75225.14
$ws.Range("J137").Value = 104257.6
$ws.Range("L137").Value = 312772.8
$ws.Range("N137").Value = -322972.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 25194.8
$ws.Range("I58").Value = 19993.5
$ws.Range("J58").Value = 46000
$ws.Range("K58").Value = 19993.5
$ws.Range("L58").Value = 46000
$ws.Range("M58").Value = -19716.5
$ws.Range("N58").Value = -46554
$ws.Range("H132").Value = 8002.8
$ws.Range("J132").Value = 13007
$ws.Range("L132").Value = 39021
$ws.Range("N132").Value = -44081

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3499.111
$ws.Range("J46").Value = 3811.5
$ws.Range("L46").Value = 3811.5
$ws.Range("N46").Value = -4187.5
$ws.Range("H55").Value = 7145727
$ws.Range("I55").Value = 10000194
$ws.Range("K55").Value = 10000194
$ws.Range("M55").Value = -10000021
$ws.Range("H131").Value = 129999
$ws.Range("J131").Value = 129999
$ws.Range("L131").Value = 129999
$ws.Range("N131").Value = -140079
$ws.Range("H136").Value = 5356.6875
$ws.Range("I136").Value = 3533.9443
$ws.Range("J136").Value = 7700.2144
$ws.Range("K136").Value = 10601.8329
$ws.Range("L136").Value = 23100.6432
$ws.Range("M136").Value = -8051.832900000001
$ws.Range("N136").Value = -28200.6432

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1500
$ws.Range("I14").Value = 1500
$ws.Range("K14").Value = 1500
$ws.Range("M14").Value = -1332
$ws.Range("H64").Value = 31109.111
$ws.Range("J64").Value = 33514
$ws.Range("L64").Value = 33514
$ws.Range("N64").Value = -34010
$ws.Range("H67").Value = 31109.111
$ws.Range("J67").Value = 33514
$ws.Range("L67").Value = 33514
$ws.Range("N67").Value = -35230
$ws.Range("H100").Value = 1042.4445
$ws.Range("I100").Value = 805.75
$ws.Range("K100").Value = 1611.5
$ws.Range("M100").Value = -1070.5
$ws.Range("H132").Value = 3774.3845
$ws.Range("I132").Value = 2296.7222
$ws.Range("K132").Value = 6890.1666
$ws.Range("M132").Value = -4360.1666

# --- Remove cells that no longer exist after the edit (ALC row 82 & 85) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M82").ClearContents()
$ws.Range("M85").ClearContents()
